# Add a new "Rights violations" field/column to the Events table.
# The Events table header lives in row 42 (columns A:W); the single
# Events data row is row 43; rows 1-41, 44 and 45 only need the new
# column's blank, formatted cell so the used range properly extends
# from column W to column X.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last column (W) into the new
# column (X) for every row of the sheet - this keeps fonts/borders
# consistent with the rest of the table without disturbing any
# existing cell.
$ws.Range("W1:W45").Copy()
$ws.Range("X1:X45").PasteSpecial(-4122)

# New header cell for the added field, in the "Events" header row.
$ws.Range("X42").Value = "Rights violations"

# New data value for the single Events row - 2 perpetrators/rights
# violations linked, matching the column's data style (numeric, not
# shared-string).
$ws.Range("X43").Value = 2

Write-Output "Added 'Rights violations' column to Events table"
